$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.311.67"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").Value = "1.903.08"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -1.48%  "
$ws.Range("D5").Value = "'315.46"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").Value = "'0.5141"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "'0.3934"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").Value = "'0.08451"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "'42.45"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.258"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.902.86"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'7.361"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "'93.21"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "'0.00001106"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'0.06733"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "'17.93"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "'6.029"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "29.282.92"
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "'2.215"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").Value = "2.120.90"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("D27").Value = "'160.31"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").Value = "'20.92"
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("D29").Value = "'2.448"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("D30").Value = "'127.74"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'1.060"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").Value = "'6.017"
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("D34").Value = "'3.645"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'0.02479"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D36").Value = "'0.06600"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "'9.132"
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("D38").Value = "'0.2198"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'1.236"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("D40").Value = "'5.127"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").Value = "'0.6515"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("D42").Value = "'1.233"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "'0.6059"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'13.25"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "'2.059"
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("D48").Value = "'1.230"
$ws.Range("D49").Value = "'123.07"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").Value = "'77.87"
$ws.Range("E51").Value = "  +1.60%  "
